$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Angpt1", "Itgb1", "ECs", 3, 1, 20.22494433333334, 60.67483300000001, 0.9541201174409912, 0.9683117878380343, 3, 1, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 2795.473306670917, 25159.25976003825, 0.2146029534793088, 0.2398120009040031),
    @("FAPs", "Angpt1", "Itgb1", "FAPs", 3, 1, 20.22494433333334, 60.67483300000001, 0.9541201174409912, 0.9683117878380343, 3, 1, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 2991.559090352297, 26924.03181317067, 0.2296560710365078, 0.2566333828221374),
    @("FAPs", "Angpt1", "Itgb1", "M1", 3, 1, 20.22494433333334, 60.67483300000001, 0.9541201174409912, 0.9683117878380343, 3, 1, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 1502.159551057635, 13519.43595951872, 0.1153178159436706, 0.1288640055179742),
    @("FAPs", "Angpt1", "Itgb1", "M2", 3, 1, 20.22494433333334, 60.67483300000001, 0.9541201174409912, 0.9683117878380343, 3, 1, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 1716.22716420458, 15446.04447784122, 0.1317513629626937, 0.1472279736213391),
    @("FAPs", "Angpt1", "Itgb1", "sCs", 3, 1, 20.22494433333334, 60.67483300000001, 0.9541201174409912, 0.9683117878380343, 2, 1, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 3423.195109564851, 20539.17065738911, 0.2627919140188104, 0.1957744249725804),
    @("M1", "Angpt1", "Itgb1", "ECs", 1, 0.3333333333333333, 0.040521, 0.121563, 0.001911594941455862, 0.00194002818046413, 3, 1, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 5.600775556792001, 50.40698001112801, 0.0004299604554956948, 0.0004804671858906201),
    @("M1", "Angpt1", "Itgb1", "FAPs", 1, 0.3333333333333333, 0.040521, 0.121563, 0.001911594941455862, 0.00194002818046413, 3, 1, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 5.993636565929999, 53.94272909337, 0.0004601196176907646, 0.0005141690940625661),
    @("M1", "Angpt1", "Itgb1", "M1", 1, 0.3333333333333333, 0.040521, 0.121563, 0.001911594941455862, 0.00194002818046413, 3, 1, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 3.009600726964, 27.086406542676, 0.0002310410917745819, 0.0002581810996131047),
    @("M1", "Angpt1", "Itgb1", "M2", 1, 0.3333333333333333, 0.040521, 0.121563, 0.001911594941455862, 0.00194002818046413, 3, 1, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 3.438488619527001, 30.94639757574301, 0.0002639659665125726, 0.0002949736039212641),
    @("M1", "Angpt1", "Itgb1", "sCs", 1, 0.3333333333333333, 0.040521, 0.121563, 0.001911594941455862, 0.00194002818046413, 2, 1, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 6.858426245755501, 41.150557474533, 0.0005265078099822482, 0.0003922371969765751),
    @("sCs", "Angpt1", "Itgb1", "ECs", 2, 1, 0.9320170000000001, 1.864034, 0.04396828761755308, 0.02974818398150157, 3, 1, 138.2190853333333, 414.657256, 0.2249223651785973, 0.2476599003709697, 128.8225372551173, 772.9352235307041, 0.009889451243792874, 0.007367432281075955),
    @("sCs", "Angpt1", "Itgb1", "FAPs", 2, 1, 0.9320170000000001, 1.864034, 0.04396828761755308, 0.02974818398150157, 3, 1, 147.91433, 443.74299, 0.2406993279341593, 0.2650317656414439, 137.85867010361, 827.15202062166, 0.01058313727996085, 0.007884213725243878),
    @("sCs", "Angpt1", "Itgb1", "M1", 2, 1, 0.9320170000000001, 1.864034, 0.04396828761755308, 0.02974818398150157, 3, 1, 74.27261733333333, 222.817852, 0.1208629960061633, 0.1330811078998542, 69.22334198916134, 415.340051934968, 0.005314138970718159, 0.003958921282266924),
    @("sCs", "Angpt1", "Itgb1", "M2", 2, 1, 0.9320170000000001, 1.864034, 0.04396828761755308, 0.02974818398150157, 3, 1, 84.85695366666668, 254.570861, 0.138086767645209, 0.1520460408212704, 79.08812338554569, 474.5287403132741, 0.006071438716002773, 0.00452309359601005),
    @("sCs", "Angpt1", "Itgb1", "sCs", 2, 1, 0.9320170000000001, 1.864034, 0.04396828761755308, 0.02974818398150157, 2, 1, 169.2560955, 338.512191, 0.275428543235871, 0.2021811852664618, 157.7495583596235, 630.9982334384941, 0.01211012140707843, 0.006014523096904758)
)

# Clear any existing data rows below the header before writing the new table
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -gt 1) {
    $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 20)).ClearContents()
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value2 = $row[$j]
    }
}
